$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-01-20 Monday" "2025-01-21 Tuesday"

Replace-Text "715÷7=102, 1" "452÷8=56, 4"
Replace-Text "645÷7=92, 1" "105÷6=17, 3"
Replace-Text "638÷6=106, 2" "148÷8=18, 4"
Replace-Text "203÷7=29, 0" "885÷8=110, 5"
Replace-Text "737÷2=368, 1" "643÷2=321, 1"

Replace-Text "419÷5=83, 4" "112÷7=16, 0"
Replace-Text "866÷4=216, 2" "584÷6=97, 2"
Replace-Text "384÷7=54, 6" "309÷6=51, 3"
Replace-Text "573÷5=114, 3" "262÷6=43, 4"
Replace-Text "653÷9=72, 5" "583÷5=116, 3"

Replace-Text "511÷9=56, 7" "857÷6=142, 5"
Replace-Text "635÷4=158, 3" "138÷4=34, 2"
Replace-Text "964÷8=120, 4" "138÷3=46, 0"
Replace-Text "355÷6=59, 1" "374÷4=93, 2"
Replace-Text "871÷3=290, 1" "215÷5=43, 0"

Replace-Text "651÷9=72, 3" "321÷4=80, 1"
Replace-Text "412÷7=58, 6" "136÷3=45, 1"
Replace-Text "664÷8=83, 0" "157÷3=52, 1"
Replace-Text "293÷4=73, 1" "738÷7=105, 3"
Replace-Text "670÷9=74, 4" "387÷6=64, 3"

Replace-Text "424÷5=84, 4" "361÷7=51, 4"
Replace-Text "323÷6=53, 5" "264÷5=52, 4"
Replace-Text "704÷5=140, 4" "422÷4=105, 2"
Replace-Text "608÷9=67, 5" "855÷6=142, 3"
Replace-Text "560÷3=186, 2" "311÷3=103, 2"
